$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("existing_stock")

$ws.Cells.Item(9,16).Value = "Aggregated Plant - IRENA Gap - way/240959264-220_Missing Hydro Capacity"
$ws.Cells.Item(10,16).Value = "Aggregated Plant - IRENA Gap - CH18-220_Missing Hydro Capacity"
$ws.Cells.Item(11,3).Value = "e_CH18-220"
$ws.Cells.Item(11,5).Value = 0.30866755954535602
$ws.Cells.Item(11,7).Value = 2750
$ws.Cells.Item(11,8).Value = 55.000000000000007
$ws.Cells.Item(11,9).Value = 2.1
$ws.Cells.Item(11,16).Value = "Aggregated Plant - IRENA Gap - way/238138373-380_Missing Hydro Capacity"
$ws.Cells.Item(12,3).Value = "e_w234983117-220"
$ws.Cells.Item(12,5).Value = 0.3629913849272427
$ws.Cells.Item(12,7).Value = 2750
$ws.Cells.Item(12,8).Value = 55.000000000000007
$ws.Cells.Item(12,9).Value = 2.1
$ws.Cells.Item(12,16).Value = "Aggregated Plant - IRENA Gap - way/234983117-220_Missing Hydro Capacity"
$ws.Cells.Item(13,3).Value = "e_r7933294-380"
$ws.Cells.Item(13,5).Value = 0.12488235719973945
$ws.Cells.Item(13,7).Value = 3162.5000000000005
$ws.Cells.Item(13,8).Value = 60.500000000000014
$ws.Cells.Item(13,9).Value = 2.3100000000000005
$ws.Cells.Item(13,16).Value = "Aggregated Plant - IRENA Gap - relation/7933294-380_Missing Hydro Capacity"
$ws.Cells.Item(14,3).Value = "e_w240959264-220"
$ws.Cells.Item(14,5).Value = 0.1117697096937668
$ws.Cells.Item(14,7).Value = 3162.5
$ws.Cells.Item(14,8).Value = 60.500000000000014
$ws.Cells.Item(14,9).Value = 2.3100000000000005
$ws.Cells.Item(15,3).Value = "e_w238138373-380"
$ws.Cells.Item(15,5).Value = 0.18732353579960917
$ws.Cells.Item(15,7).Value = 3162.5000000000009
$ws.Cells.Item(136,16).Value = "Aggregated Plant - IRENA Gap - CHE_18_Missing Solar Capacity"
$ws.Cells.Item(137,16).Value = "Aggregated Plant - IRENA Gap - CHE_2_Missing Solar Capacity"
$ws.Cells.Item(138,16).Value = "Aggregated Plant - IRENA Gap - CHE_11_Missing Solar Capacity"
$ws.Cells.Item(139,16).Value = "Aggregated Plant - IRENA Gap - CHE_0_Missing Solar Capacity"
$ws.Cells.Item(140,16).Value = "Aggregated Plant - IRENA Gap - CHE_23_Missing Solar Capacity"
$ws.Cells.Item(141,16).Value = "Aggregated Plant - IRENA Gap - CHE_6_Missing Solar Capacity"
$ws.Cells.Item(142,16).Value = "Aggregated Plant - IRENA Gap - CHE_13_Missing Solar Capacity"
$ws.Cells.Item(143,16).Value = "Aggregated Plant - IRENA Gap - CHE_15_Missing Solar Capacity"
$ws.Cells.Item(144,16).Value = "Aggregated Plant - IRENA Gap - CHE_7_Missing Solar Capacity"
$ws.Cells.Item(145,16).Value = "Aggregated Plant - IRENA Gap - CHE_12_Missing Solar Capacity"
$ws.Cells.Item(146,16).Value = "Aggregated Plant - IRENA Gap - CHE_25_Missing Solar Capacity"
$ws.Cells.Item(147,16).Value = "Aggregated Plant - IRENA Gap - CHE_19_Missing Solar Capacity"
$ws.Cells.Item(148,16).Value = "Aggregated Plant - IRENA Gap - CHE_10_Missing Solar Capacity"
$ws.Cells.Item(150,3).Value = "elc_spv-CHE_0017"
$ws.Cells.Item(150,5).Value = 0.15226887751132734
$ws.Cells.Item(150,16).Value = "Aggregated Plant - IRENA Gap - CHE_21_Missing Solar Capacity"
$ws.Cells.Item(151,3).Value = "elc_spv-CHE_0013"
$ws.Cells.Item(151,5).Value = 0.17206733071733712
$ws.Cells.Item(151,16).Value = "Aggregated Plant - IRENA Gap - CHE_22_Missing Solar Capacity"
$ws.Cells.Item(152,3).Value = "elc_spv-CHE_0005"
$ws.Cells.Item(152,5).Value = 0.20006982412215921
$ws.Cells.Item(152,16).Value = "Aggregated Plant - IRENA Gap - CHE_9_Missing Solar Capacity"
$ws.Cells.Item(153,3).Value = "elc_spv-CHE_0012"
$ws.Cells.Item(153,5).Value = 0.13549669849969209
$ws.Cells.Item(153,16).Value = "Aggregated Plant - IRENA Gap - CHE_20_Missing Solar Capacity"
$ws.Cells.Item(154,3).Value = "elc_spv-CHE_0014"
$ws.Cells.Item(154,5).Value = 0.18231505170803797
$ws.Cells.Item(154,16).Value = "Aggregated Plant - IRENA Gap - CHE_3_Missing Solar Capacity"
$ws.Cells.Item(155,3).Value = "elc_spv-CHE_0001"
$ws.Cells.Item(155,5).Value = 0.20988535532947597
$ws.Cells.Item(155,16).Value = "Aggregated Plant - IRENA Gap - CHE_1_Missing Solar Capacity"
$ws.Cells.Item(156,3).Value = "elc_spv-CHE_0015"
$ws.Cells.Item(156,5).Value = 0.1701049810444224
$ws.Cells.Item(156,16).Value = "Aggregated Plant - IRENA Gap - CHE_24_Missing Solar Capacity"
$ws.Cells.Item(157,3).Value = "elc_spv-CHE_0003"
$ws.Cells.Item(157,5).Value = 0.16085025627375071
$ws.Cells.Item(157,16).Value = "Aggregated Plant - IRENA Gap - CHE_5_Missing Solar Capacity"
$ws.Cells.Item(158,3).Value = "elc_spv-CHE_0009"
$ws.Cells.Item(158,5).Value = 0.19532613932247714
$ws.Cells.Item(158,16).Value = "Aggregated Plant - IRENA Gap - CHE_4_Missing Solar Capacity"
$ws.Cells.Item(159,3).Value = "elc_spv-CHE_0011"
$ws.Cells.Item(159,5).Value = 0.16209575724687297
$ws.Cells.Item(159,16).Value = "Aggregated Plant - IRENA Gap - CHE_8_Missing Solar Capacity"
$ws.Cells.Item(160,3).Value = "elc_spv-CHE_0004"
$ws.Cells.Item(160,5).Value = 0.19745398836539674
$ws.Cells.Item(160,16).Value = "Aggregated Plant - IRENA Gap - CHE_17_Missing Solar Capacity"
$ws.Cells.Item(161,3).Value = "elc_spv-CHE_0021"
$ws.Cells.Item(161,5).Value = 0.15273795001145538
$ws.Cells.Item(162,3).Value = "elc_spv-CHE_0006"
$ws.Cells.Item(162,5).Value = 0.21640319337561012
$ws.Cells.Item(163,3).Value = "elc_spv-CHE_0000"
$ws.Cells.Item(163,5).Value = 0.19247860444770779
$ws.Cells.Item(164,3).Value = "elc_spv-CHE_0010"
$ws.Cells.Item(164,5).Value = 0.19228757088918788
$ws.Cells.Item(165,3).Value = "elc_spv-CHE_0018"
$ws.Cells.Item(165,5).Value = 0.19699531868281184
$ws.Cells.Item(166,3).Value = "elc_spv-CHE_0002"
$ws.Cells.Item(166,5).Value = 0.21063530390326943
$ws.Cells.Item(167,3).Value = "elc_spv-CHE_0019"
$ws.Cells.Item(167,5).Value = 0.16485344960649678
$ws.Cells.Item(168,3).Value = "elc_spv-CHE_0020"
$ws.Cells.Item(168,5).Value = 0.15456128021356608
$ws.Cells.Item(169,3).Value = "elc_spv-CHE_0007"
$ws.Cells.Item(169,5).Value = 0.16629376698088194
$ws.Cells.Item(170,3).Value = "elc_spv-CHE_0023"
$ws.Cells.Item(170,5).Value = 0.1579516530219513
$ws.Cells.Item(171,3).Value = "elc_spv-CHE_0024"
$ws.Cells.Item(171,5).Value = 0.16568094645652107
$ws.Cells.Item(172,3).Value = "elc_spv-CHE_0008"
$ws.Cells.Item(172,5).Value = 0.1534481787364477
$ws.Cells.Item(173,3).Value = "elc_spv-CHE_0022"
$ws.Cells.Item(173,5).Value = 0.21381383751804844
$ws.Cells.Item(174,3).Value = "elc_spv-CHE_0025"
$ws.Cells.Item(174,5).Value = 0.13652468601509371
